$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the y_fixStart / y_nrSteps / alienID values on trial row 3 (sheet row 3)
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Leave the selection on the last-edited cell, matching the saved view state
$ws.Range("E3").Select()
